# Re-apply the word list from a re-processing/re-ordering run of the source script.
# Only the shared string ordering changed (same multiset of goods words);
# rows/counts in column B stay fixed - the word shown in column A for a given
# row is updated to match the new ordering produced by the rerun.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "Крымскую соль"
$ws.Range("A17").Value = "полотно"
$ws.Range("A18").Value = "колеса"
$ws.Range("A19").Value = "сено"
$ws.Range("A21").Value = "парча"
$ws.Range("A22").Value = "табак"
$ws.Range("A23").Value = "позумент"
$ws.Range("A24").Value = "шелк"
$ws.Range("A25").Value = "чулок"
$ws.Range("A26").Value = "выбойка"
$ws.Range("A27").Value = "сахар"
$ws.Range("A28").Value = "лыко"
$ws.Range("A29").Value = "лес"
$ws.Range("A30").Value = "коса"
$ws.Range("A31").Value = "ладан"
$ws.Range("A32").Value = "сапог"
$ws.Range("A33").Value = "китайка"
$ws.Range("A34").Value = "сани"
$ws.Range("A35").Value = "гвоздь"
$ws.Range("A36").Value = "конь"
$ws.Range("A37").Value = "горшок"
$ws.Range("A38").Value = "веревка"
$ws.Range("A39").Value = "платок"
$ws.Range("A40").Value = "рогожа"
$ws.Range("A41").Value = "замок"
$ws.Range("A42").Value = "овца"
$ws.Range("A43").Value = "обод"
$ws.Range("A45").Value = "дуга"
$ws.Range("A46").Value = "бечева"
$ws.Range("A47").Value = "котел"
$ws.Range("A48").Value = "гумми"
$ws.Range("A49").Value = "сковорода"
$ws.Range("A50").Value = "брусья"
$ws.Range("A51").Value = "роза"
$ws.Range("A52").Value = "сосуд"
$ws.Range("A54").Value = "скотский кожа"
$ws.Range("A55").Value = "покроми"
$ws.Range("A56").Value = "нитка"
